$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "26.697.11"
Set-TextValue $ws.Range("E2") "  +3.69%  "
Set-TextValue $ws.Range("D3") "1.867.01"
Set-TextValue $ws.Range("E3") "  +2.79%  "
Set-TextValue $ws.Range("D4") "1.002"
Set-TextValue $ws.Range("E4") "  +0.19%  "
Set-TextValue $ws.Range("D5") "276.89"
Set-TextValue $ws.Range("E5") "  -0.64%  "
Set-TextValue $ws.Range("D6") "1.002"
Set-TextValue $ws.Range("E6") "  +0.16%  "
Set-TextValue $ws.Range("D7") "0.5276"
Set-TextValue $ws.Range("E7") "  +3.67%  "
Set-TextValue $ws.Range("D8") "0.3412"
Set-TextValue $ws.Range("E8") "  -3.53%  "
Set-TextValue $ws.Range("D9") "0.06924"
Set-TextValue $ws.Range("E9") "  +3.90%  "
Set-TextValue $ws.Range("D10") "19.98"
Set-TextValue $ws.Range("E10") "  -0.38%  "
Set-TextValue $ws.Range("D11") "0.8011"
Set-TextValue $ws.Range("E11") "  -3.08%  "
Set-TextValue $ws.Range("D12") "0.07738"
Set-TextValue $ws.Range("E12") "  -2.08%  "
Set-TextValue $ws.Range("D13") "1.884.14"
Set-TextValue $ws.Range("E13") "  +3.80%  "
Set-TextValue $ws.Range("D14") "89.95"
Set-TextValue $ws.Range("E14") "  +2.52%  "
Set-TextValue $ws.Range("D15") "5.159"
Set-TextValue $ws.Range("E15") "  +1.66%  "
Set-TextValue $ws.Range("B16") "BinanceUSD"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D16") "1.002"
Set-TextValue $ws.Range("E16") "  +0.19%  "
Set-TextValue $ws.Range("B17") "Avalanche"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D17") "14.52"
Set-TextValue $ws.Range("E17") "  +3.16%  "
Set-TextValue $ws.Range("D18") "0.000008024"
Set-TextValue $ws.Range("E18") "  -0.12%  "
Set-TextValue $ws.Range("E19") "  +0.12%  "
Set-TextValue $ws.Range("D20") "26.745.32"
Set-TextValue $ws.Range("E20") "  +3.72%  "
Set-TextValue $ws.Range("D21") "2.123.81"
Set-TextValue $ws.Range("E21") "  +4.43%  "
Set-TextValue $ws.Range("D22") "4.745"
Set-TextValue $ws.Range("E22") "  -0.05%  "
Set-TextValue $ws.Range("D23") "10.01"
Set-TextValue $ws.Range("E23") "  +0.24%  "
Set-TextValue $ws.Range("D24") "6.170"
Set-TextValue $ws.Range("E24") "  +0.89%  "
Set-TextValue $ws.Range("D25") "2.358"
Set-TextValue $ws.Range("E25") "  +5.60%  "
Set-TextValue $ws.Range("D26") "146.02"
Set-TextValue $ws.Range("E26") "  +2.64%  "
Set-TextValue $ws.Range("D27") "17.30"
Set-TextValue $ws.Range("E27") "  +1.13%  "
Set-TextValue $ws.Range("D28") "1.652"
Set-TextValue $ws.Range("E28") "  -0.87%  "
Set-TextValue $ws.Range("D29") "113.00"
Set-TextValue $ws.Range("E29") "  +3.37%  "
Set-TextValue $ws.Range("D30") "4.326"
Set-TextValue $ws.Range("E30") "  -0.02%  "
Set-TextValue $ws.Range("D31") "4.333"
Set-TextValue $ws.Range("E31") "  +2.13%  "
Set-TextValue $ws.Range("D32") "0.08889"
Set-TextValue $ws.Range("E32") "  +1.33%  "
Set-TextValue $ws.Range("D33") "0.04933"
Set-TextValue $ws.Range("E33") "  +0.52%  "
Set-TextValue $ws.Range("D34") "1.161"
Set-TextValue $ws.Range("E34") "  +1.92%  "
Set-TextValue $ws.Range("D35") "0.7277"
Set-TextValue $ws.Range("E35") "  -0.39%  "
Set-TextValue $ws.Range("D36") "2.889"
Set-TextValue $ws.Range("E36") "  +0.57%  "
Set-TextValue $ws.Range("D37") "3.255"
Set-TextValue $ws.Range("E37") "  +3.54%  "
Set-TextValue $ws.Range("D38") "0.01851"
Set-TextValue $ws.Range("E38") "  -0.07%  "
Set-TextValue $ws.Range("D39") "2.311"
Set-TextValue $ws.Range("E39") "  -3.41%  "
Set-TextValue $ws.Range("D40") "0.5130"
Set-TextValue $ws.Range("E40") "  -0.30%  "
Set-TextValue $ws.Range("D41") "0.9458"
Set-TextValue $ws.Range("E41") "  -1.95%  "
Set-TextValue $ws.Range("D42") "116.20"
Set-TextValue $ws.Range("E42") "  +4.55%  "
Set-TextValue $ws.Range("D43") "6.139"
Set-TextValue $ws.Range("E43") "  -1.29%  "
Set-TextValue $ws.Range("D44") "8.082"
Set-TextValue $ws.Range("E44") "  +0.59%  "
Set-TextValue $ws.Range("E45") "  +0.10%  "
Set-TextValue $ws.Range("D46") "0.4449"
Set-TextValue $ws.Range("E46") "  -2.42%  "
Set-TextValue $ws.Range("D47") "0.1338"
Set-TextValue $ws.Range("E47") "  -2.24%  "
Set-TextValue $ws.Range("B48") "EnergySwap"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D48") "9.290"
Set-TextValue $ws.Range("E48") "  +1.03%  "
Set-TextValue $ws.Range("B49") "Cronos"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D49") "0.06019"
Set-TextValue $ws.Range("E49") "  +3.21%  "
Set-TextValue $ws.Range("B50") "Elrond"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue $ws.Range("D50") "36.29"
Set-TextValue $ws.Range("E50") "  -0.77%  "
Set-TextValue $ws.Range("D51") "1.484"
Set-TextValue $ws.Range("E51") "  -1.22%  "
